$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value as genuine TEXT (shared-string) while
# keeping the destination cell's existing style untouched. A direct
# `.Value = "123"` assignment auto-coerces digit-only strings to a Number,
# which would drop the `t="s"` shared-string typing. Instead, build the text
# in a scratch cell via a formula that evaluates to a string (so no
# NumberFormat change is needed -- that would mint an unused extra cell
# style), copy it, and paste-special just the value into the destination so
# the destination keeps its own existing style index.
function Set-TextValue {
    param($cell, [string]$text)

    $scratch = $ws.Cells.Item(1048576, 16384)
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# Row 2: 20113120 / AICE CHOCO ALMOND 90 / .. / .. / 15 / ..
#     -> 20032368 / WALL'S MAGNUM ALMD80  / .. / .. / 1  / ..
Set-TextValue $ws.Cells.Item(2, 1) "20032368"
$ws.Cells.Item(2, 2).Value = "WALL'S MAGNUM ALMD80"
Set-TextValue $ws.Cells.Item(2, 5) "1"

# Row 3: 20134511 / AICE CLSC CHO ALMD90 / .. / .. / 91 / RT,(E-1B)
#     -> 20032366 / WALL'S MAGNUM CLAS80 / .. / .. / 2  / RT,(E-3B)
Set-TextValue $ws.Cells.Item(3, 1) "20032366"
$ws.Cells.Item(3, 2).Value = "WALL'S MAGNUM CLAS80"
Set-TextValue $ws.Cells.Item(3, 5) "2"
$ws.Cells.Item(3, 6).Value = "RT,(E-3B)"

# Columns D and E both end up width 3 (was 3 and 4).
$ws.Columns("D:E").ColumnWidth = 2.1666666666666665
